$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 9.380719000000001
$ws.Range("H2").Value = 28.142157
$ws.Range("I2").Value = 0.03679977590837273
$ws.Range("J2").Value = 0.03679977590837273
$ws.Range("O2").Value = 0.01611173663836548
$ws.Range("P2").Value = 0.01611173663836548
$ws.Range("Q2").Value = 0.597561181019
$ws.Range("R2").Value = 5.378050629171
$ws.Range("S2").Value = 0.0005929082977865683
$ws.Range("T2").Value = 0.0005929082977865681
$ws.Range("G3").Value = 9.380719000000001
$ws.Range("H3").Value = 28.142157
$ws.Range("I3").Value = 0.03679977590837273
$ws.Range("J3").Value = 0.03679977590837273
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 2.701496333333333
$ws.Range("N3").Value = 8.104489000000001
$ws.Range("O3").Value = 0.68328279700753
$ws.Range("P3").Value = 0.68328279700753
$ws.Range("Q3").Value = 25.34197798253034
$ws.Range("R3").Value = 228.077801842773
$ws.Range("S3").Value = 0.02514465381192324
$ws.Range("T3").Value = 0.02514465381192324
$ws.Range("G4").Value = 9.380719000000001
$ws.Range("H4").Value = 28.142157
$ws.Range("I4").Value = 0.03679977590837273
$ws.Range("J4").Value = 0.03679977590837273
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.188504333333333
$ws.Range("N4").Value = 3.565513
$ws.Range("O4").Value = 0.3006054663541045
$ws.Range("P4").Value = 0.3006054663541044
$ws.Range("Q4").Value = 11.14902518128233
$ws.Range("R4").Value = 100.341226631541
$ws.Range("S4").Value = 0.01106221379866292
$ws.Range("T4").Value = 0.01106221379866292
$ws.Range("I5").Value = 0.3547860986448385
$ws.Range("J5").Value = 0.3547860986448385
$ws.Range("O5").Value = 0.01611173663836548
$ws.Range("P5").Value = 0.01611173663836548
$ws.Range("S5").Value = 0.005716220184318794
$ws.Range("T5").Value = 0.005716220184318793
$ws.Range("I6").Value = 0.3547860986448385
$ws.Range("J6").Value = 0.3547860986448385
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.701496333333333
$ws.Range("N6").Value = 8.104489000000001
$ws.Range("O6").Value = 0.68328279700753
$ws.Range("P6").Value = 0.68328279700753
$ws.Range("Q6").Value = 244.321637249962
$ws.Range("R6").Value = 2198.894735249658
$ws.Range("S6").Value = 0.2424192378214347
$ws.Range("T6").Value = 0.2424192378214347
$ws.Range("I7").Value = 0.3547860986448385
$ws.Range("J7").Value = 0.3547860986448385
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.188504333333333
$ws.Range("N7").Value = 3.565513
$ws.Range("O7").Value = 0.3006054663541045
$ws.Range("P7").Value = 0.3006054663541044
$ws.Range("Q7").Value = 107.487587902954
$ws.Range("R7").Value = 967.3882911265861
$ws.Range("S7").Value = 0.106650640639085
$ws.Range("T7").Value = 0.106650640639085
$ws.Range("G8").Value = 100.179423
$ws.Range("H8").Value = 300.538269
$ws.Range("I8").Value = 0.3929954960840508
$ws.Range("J8").Value = 0.3929954960840508
$ws.Range("O8").Value = 0.01611173663836548
$ws.Range("P8").Value = 0.01611173663836548
$ws.Range("Q8").Value = 6.381529424522999
$ws.Range("R8").Value = 57.433764820707
$ws.Range("S8").Value = 0.006331839932970019
$ws.Range("T8").Value = 0.006331839932970017
$ws.Range("G9").Value = 100.179423
$ws.Range("H9").Value = 300.538269
$ws.Range("I9").Value = 0.3929954960840508
$ws.Range("J9").Value = 0.3929954960840508
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 2.701496333333333
$ws.Range("N9").Value = 8.104489000000001
$ws.Range("O9").Value = 0.68328279700753
$ws.Range("P9").Value = 0.68328279700753
$ws.Range("Q9").Value = 270.634343909949
$ws.Range("R9").Value = 2435.709095189542
$ws.Range("S9").Value = 0.268527061775672
$ws.Range("T9").Value = 0.268527061775672
$ws.Range("G10").Value = 100.179423
$ws.Range("H10").Value = 300.538269
$ws.Range("I10").Value = 0.3929954960840508
$ws.Range("J10").Value = 0.3929954960840508
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.188504333333333
$ws.Range("N10").Value = 3.565513
$ws.Range("O10").Value = 0.3006054663541045
$ws.Range("P10").Value = 0.3006054663541044
$ws.Range("Q10").Value = 119.063678346333
$ws.Range("R10").Value = 1071.573105116997
$ws.Range("S10").Value = 0.1181365943754087
$ws.Range("T10").Value = 0.1181365943754087
$ws.Range("G11").Value = 1.427630666666667
$ws.Range("H11").Value = 4.282892
$ws.Range("I11").Value = 0.005600475679236752
$ws.Range("J11").Value = 0.005600475679236752
$ws.Range("O11").Value = 0.01611173663836548
$ws.Range("P11").Value = 0.01611173663836548
$ws.Range("Q11").Value = 0.09094150109733333
$ws.Range("R11").Value = 0.8184735098760001
$ws.Range("S11").Value = 0.00009023338919343357
$ws.Range("T11").Value = 0.00009023338919343355
$ws.Range("G12").Value = 1.427630666666667
$ws.Range("H12").Value = 4.282892
$ws.Range("I12").Value = 0.005600475679236752
$ws.Range("J12").Value = 0.005600475679236752
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 2.701496333333333
$ws.Range("N12").Value = 8.104489000000001
$ws.Range("O12").Value = 0.68328279700753
$ws.Range("P12").Value = 0.68328279700753
$ws.Range("Q12").Value = 3.856739011354223
$ws.Range("R12").Value = 34.71065110218801
$ws.Range("S12").Value = 0.003826708686681534
$ws.Range("T12").Value = 0.003826708686681534
$ws.Range("G13").Value = 1.427630666666667
$ws.Range("H13").Value = 4.282892
$ws.Range("I13").Value = 0.005600475679236752
$ws.Range("J13").Value = 0.005600475679236752
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 1.188504333333333
$ws.Range("N13").Value = 3.565513
$ws.Range("O13").Value = 0.3006054663541045
$ws.Range("P13").Value = 0.3006054663541044
$ws.Range("Q13").Value = 1.696745233732889
$ws.Range("R13").Value = 15.270707103596
$ws.Range("S13").Value = 0.001683533603361784
$ws.Range("T13").Value = 0.001683533603361783
$ws.Range("G14").Value = 53.48524799999999
$ws.Range("H14").Value = 160.455744
$ws.Range("I14").Value = 0.2098181536835013
$ws.Range("J14").Value = 0.2098181536835013
$ws.Range("O14").Value = 0.01611173663836548
$ws.Range("P14").Value = 0.01611173663836548
$ws.Range("Q14").Value = 3.407063782847999
$ws.Range("R14").Value = 30.663574045632
$ws.Range("S14").Value = 0.003380534834096666
$ws.Range("T14").Value = 0.003380534834096666
$ws.Range("G15").Value = 53.48524799999999
$ws.Range("H15").Value = 160.455744
$ws.Range("I15").Value = 0.2098181536835013
$ws.Range("J15").Value = 0.2098181536835013
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 2.701496333333333
$ws.Range("N15").Value = 8.104489000000001
$ws.Range("O15").Value = 0.68328279700753
$ws.Range("P15").Value = 0.68328279700753
$ws.Range("Q15").Value = 144.490201359424
$ws.Range("R15").Value = 1300.411812234816
$ws.Range("S15").Value = 0.1433651349118185
$ws.Range("T15").Value = 0.1433651349118185
$ws.Range("G16").Value = 53.48524799999999
$ws.Range("H16").Value = 160.455744
$ws.Range("I16").Value = 0.2098181536835013
$ws.Range("J16").Value = 0.2098181536835013
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 1.188504333333333
$ws.Range("N16").Value = 3.565513
$ws.Range("O16").Value = 0.3006054663541045
$ws.Range("P16").Value = 0.3006054663541044
$ws.Range("Q16").Value = 63.56744901740799
$ws.Range("R16").Value = 572.1070411566719
$ws.Range("S16").Value = 0.06307248393758606
$ws.Range("T16").Value = 0.06307248393758605
